$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# WEEK 3 block (rows 14-16): fill in the previously-empty "volume" column (D)
$ws.Range("D14").Value = "2x20"
$ws.Range("D15").Value = "2x20, 1x4.5, 1x2.5"
$ws.Range("D16").Value = "2x20, 1x15.6"

# WEEK 4 block (rows 20-22): fill in the previously-empty "volume" column (D)
$ws.Range("D20").Value = "1x20"
$ws.Range("D21").Value = "1x20, 1x4.5"
$ws.Range("D22").Value = "1x20, 1x11.3"

# Column D widened to fit the newly-added longer text
$ws.Columns.Item(4).ColumnWidth = 14.3

# Update the window view/selection to focus on the newly-entered D14:D16 block
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("D14:D16").Select()
